# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a few cell clears) to the FFXIV leve-profit sheets
# as described by the commit "chore: update Sheets via scheduled runner".

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @("H12", 366.66666),
    @("I12", 300),
    @("K12", 300),
    @("M12", -130),
    @("H28", 938.875),
    @("I28", 814.5862),
    @("J28", 1266.5454),
    @("K28", 814.5862),
    @("L28", 1266.5454),
    @("M28", -329.5862),
    @("N28", -2236.5454),
    @("H33", 727.8570999999999),
    @("I33", 760.7692),
    @("K33", 760.7692),
    @("M33", -531.7692),
    @("H53", 1148),
    @("I53", 160.66667),
    @("J53", 1740.4),
    @("K53", 160.66667),
    @("L53", 1740.4),
    @("M53", 476.33333),
    @("N53", -3014.4),
    @("H70", 4874.5),
    @("I70", 4750),
    @("J70", 4899.4),
    @("K70", 14250),
    @("L70", 14698.2),
    @("M70", -13980),
    @("N70", -15238.2),
    @("H73", 4874.5),
    @("I73", 4750),
    @("J73", 4899.4),
    @("K73", 14250),
    @("L73", 14698.2),
    @("M73", -13314),
    @("N73", -16570.2),
    @("H82", 6916.1665),
    @("I82", 4299.4),
    @("K82", 12898.2),
    @("M82", -12492.2),
    @("H85", 6916.1665),
    @("I85", 4299.4),
    @("K85", 12898.2),
    @("M85", -11494.2),
    @("H96", 2250),
    @("I96", 2500),
    @("J96", 2000),
    @("K96", 7500),
    @("L96", 6000),
    @("M96", -6127),
    @("N96", -8746),
    @("H132", 34177.824),
    @("I132", 2585.3333),
    @("K132", 7755.999899999999),
    @("M132", -5225.999899999999),
    @("H137", 3914.7354),
    @("I137", 2244.25),
    @("J137", 6301.143),
    @("K137", 6732.75),
    @("L137", 18903.429),
    @("M137", -4182.75),
    @("N137", -24003.429),
    @("H141", 2193.875),
    @("I141", 2193.875),
    @("K141", 6581.625),
    @("M141", -1401.625)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @("H32", 6210.8237),
    @("I32", 5776.7236),
    @("K32", 5776.7236),
    @("M32", -5489.7236),
    @("H34", 29997.666),
    @("I34", 0),
    @("J34", 29997.666),
    @("K34", 0),
    @("L34", 29997.666),
    @("N34", -30539.666),
    @("H45", 2836.625),
    @("I45", 2517.5),
    @("J45", 2943),
    @("K45", 2517.5),
    @("L45", 2943),
    @("M45", -2140.5),
    @("N45", -3697),
    @("H61", 3463.1714),
    @("I61", 2373.52),
    @("K61", 2373.52),
    @("M61", -2161.52),
    @("H74", 782.3200000000001),
    @("I74", 770.8182),
    @("K74", 770.8182),
    @("M74", 103.1818),
    @("H77", 782.3200000000001),
    @("I77", 770.8182),
    @("K77", 3854.091),
    @("M77", 513.9089999999997),
    @("H97", 3832258.2),
    @("I97", 927.8095),
    @("K97", 927.8095),
    @("M97", -431.8095),
    @("H136", 3463.1714),
    @("I136", 2373.52),
    @("K136", 7120.559999999999),
    @("M136", -4570.559999999999)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$clears = @("M34")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @("H20", 3071.8572),
    @("I20", 2442.4119),
    @("K20", 2442.4119),
    @("M20", -2195.4119),
    @("H99", 30801.857),
    @("I99", 42151.5),
    @("K99", 42151.5),
    @("M99", -40653.5),
    @("H105", 4634.9),
    @("I105", 4483.222),
    @("K105", 4483.222),
    @("M105", -2736.222)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @("H31", 4716.7646),
    @("I31", 2001.8334),
    @("J31", 6197.636),
    @("K31", 2001.8334),
    @("L31", 6197.636),
    @("M31", -1706.8334),
    @("N31", -6787.636),
    @("H34", 4716.7646),
    @("I34", 2001.8334),
    @("J34", 6197.636),
    @("K34", 2001.8334),
    @("L34", 6197.636),
    @("M34", -1799.8334),
    @("N34", -6601.636),
    @("H86", 29113.785),
    @("I86", 47336.145),
    @("K86", 47336.145),
    @("M86", -46213.145),
    @("H89", 29113.785),
    @("I89", 47336.145),
    @("K89", 236680.725),
    @("M89", -231064.725),
    @("H122", 490822.56),
    @("I122", 1135708.8),
    @("J122", 7157.9165),
    @("K122", 3407126.4),
    @("L122", 21473.7495),
    @("M122", -3404676.4),
    @("N122", -26373.7495),
    @("H132", 3113.862),
    @("I132", 2819.3462),
    @("K132", 8458.0386),
    @("M132", -5928.0386),
    @("H134", 4343.75),
    @("I134", 3058.4119),
    @("K134", 9175.235700000001),
    @("M134", -6640.235700000001)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @("H54", 145783),
    @("J54", 145783),
    @("L54", 437349),
    @("N54", -438467),
    @("H60", 303.4),
    @("I60", 275.625),
    @("K60", 826.875),
    @("M60", -575.875),
    @("H131", 5100.96),
    @("I131", 2889.7778),
    @("J131", 6344.75),
    @("K131", 8669.3334),
    @("L131", 19034.25),
    @("M131", -3629.3334),
    @("N131", -29114.25)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @("H14", 2322500),
    @("J14", 895000),
    @("L14", 895000),
    @("N14", -895336),
    @("H70", 16493),
    @("I70", 15206.429),
    @("J70", 20996),
    @("K70", 15206.429),
    @("L70", 20996),
    @("M70", -14936.429),
    @("N70", -21536),
    @("H73", 16493),
    @("I73", 15206.429),
    @("J73", 20996),
    @("K73", 15206.429),
    @("L73", 20996),
    @("M73", -14270.429),
    @("N73", -22868),
    @("H132", 1247.9375),
    @("I132", 1036.1538),
    @("J132", 2165.6667),
    @("K132", 3108.4614),
    @("L132", 6497.000100000001),
    @("M132", -578.4614000000001),
    @("N132", -11557.0001)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @("H38", 0),
    @("J38", 0),
    @("L38", 0),
    @("H82", 3687.4333),
    @("I82", 2195.55),
    @("J82", 6671.2),
    @("K82", 2195.55),
    @("L82", 6671.2),
    @("M82", -1834.55),
    @("N82", -7393.2),
    @("H85", 3687.4333),
    @("I85", 2195.55),
    @("J85", 6671.2),
    @("K85", 2195.55),
    @("L85", 6671.2),
    @("M85", -947.5500000000002),
    @("N85", -9167.200000000001),
    @("H103", 14999.5),
    @("J103", 14999.5),
    @("L103", 14999.5),
    @("N103", -17343.5),
    @("H122", 5504),
    @("I122", 2124.5),
    @("J122", 7006),
    @("K122", 6373.5),
    @("L122", 21018),
    @("M122", -3923.5),
    @("N122", -25918),
    @("H133", 78661.75),
    @("J133", 78661.75),
    @("L133", 78661.75),
    @("N133", -83721.75)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$clears = @("N38")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @("H20", 0),
    @("J20", 0),
    @("L20", 0),
    @("H136", 3816.524),
    @("I136", 2214.7273),
    @("K136", 6644.1819),
    @("M136", -4094.1819)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

$clears = @("N20")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}
